$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.217.06'
$ws.Range("E2").Value = '  -7.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.673.04'
$ws.Range("E3").Value = '  -4.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.58'
$ws.Range("E5").Value = '  -4.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5099'
$ws.Range("E6").Value = '  -12.49%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2654'
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.80'
$ws.Range("E9").Value = '  -4.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06319'
$ws.Range("E10").Value = '  -4.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07369'
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.679.08'
$ws.Range("E12").Value = '  -5.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.553'
$ws.Range("E13").Value = '  -3.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5753'
$ws.Range("E14").Value = '  -4.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.906.88'
$ws.Range("E15").Value = '  -4.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008515'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.94'
$ws.Range("E17").Value = '  -12.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.306.17'
$ws.Range("E18").Value = '  -7.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.001'
$ws.Range("E19").Value = '  -5.94%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.87'
$ws.Range("E21").Value = '  -3.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '186.64'
$ws.Range("E22").Value = '  -9.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.211'
$ws.Range("E23").Value = '  -7.07%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.67'
$ws.Range("E25").Value = '  -4.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.563'
$ws.Range("E26").Value = '  -5.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1175'
$ws.Range("E27").Value = '  -5.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.68'
$ws.Range("E28").Value = '  -2.75%  '
$ws.Range("E29").Value = '  -5.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05797'
$ws.Range("E30").Value = '  -5.43%  '
$ws.Range("E31").Value = '  -5.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.509'
$ws.Range("E32").Value = '  -5.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.502'
$ws.Range("E33").Value = '  -5.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.660'
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.001'
$ws.Range("E35").Value = '  -2.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.5977'
$ws.Range("E36").Value = '  -5.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.368'
$ws.Range("E37").Value = '  -2.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.661'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.097.25'
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01601'
$ws.Range("E40").Value = '  -4.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.911'
$ws.Range("E41").Value = '  -6.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8620'
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.36'
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.825.68'
$ws.Range("E45").Value = '  -4.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000115'
$ws.Range("E46").Value = '  +4.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.24'
$ws.Range("E47").Value = '  -4.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.066'
$ws.Range("E49").Value = '  -2.83%  '
$ws.Range("E50").Value = '  -3.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05213'
$ws.Range("E51").Value = '  -3.55%  '
